# This workbook tracks weekly Albahaca (basil) price records for the
# "Vega Central Mapocho de Santiago" market. A new weekly entry (two rows,
# "Primera" and "Segunda"-style quality records for the same date) is
# inserted at the top of the historical block that starts at row 523,
# pushing all the existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 523; this shifts the old
# rows 523:627 down to 525:629 (dimension grows from R627 to R629).
$ws.Rows("523:524").Insert()

# --- Row 523: new "Primera" record -----------------------------------
$ws.Range("A523").Value2 = 9
$ws.Range("B523").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C523").Value2 = "Metropolitana"
$ws.Range("D523").Value2 = 45209
$ws.Range("E523").Value2 = 13
$ws.Range("F523").Value2 = 100112052
$ws.Range("G523").Value2 = "Albahaca"
$ws.Range("H523").Value2 = "Sin especificar"
$ws.Range("I523").Value2 = "Primera"
$ws.Range("J523").Value2 = 160
$ws.Range("K523").Value2 = 5000
$ws.Range("L523").Value2 = 5000
$ws.Range("M523").Value2 = 5000
$ws.Range("N523").Value2 = "`$/docena de matas"
$ws.Range("O523").Value2 = "Provincia de Chacabuco"
$ws.Range("P523").Value2 = 833
$ws.Range("Q523").Value2 = 6
$ws.Range("R523").Value2 = "Hortaliza"

# --- Row 524: new "Primera" record -----------------------------------
$ws.Range("A524").Value2 = 9
$ws.Range("B524").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C524").Value2 = "Metropolitana"
$ws.Range("D524").Value2 = 45209
$ws.Range("E524").Value2 = 13
$ws.Range("F524").Value2 = 100112052
$ws.Range("G524").Value2 = "Albahaca"
$ws.Range("H524").Value2 = "Sin especificar"
$ws.Range("I524").Value2 = "Primera"
$ws.Range("J524").Value2 = 340
$ws.Range("K524").Value2 = 5000
$ws.Range("L524").Value2 = 5500
$ws.Range("M524").Value2 = 5250
$ws.Range("N524").Value2 = "`$/paquete"
$ws.Range("O524").Value2 = "Región de Arica y Parinacota"
$ws.Range("P524").Value2 = 5250
$ws.Range("Q524").Value2 = 1
$ws.Range("R524").Value2 = "Hortaliza"
